$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 95, shifting existing rows 95-100 down to 96-101.
$ws.Rows.Item(95).Insert()

# Fill in the new row 95 with the new weekly record (same pattern as neighboring rows,
# with updated date / volume / price figures).
$ws.Cells.Item(95, 1).Value = 2
$ws.Cells.Item(95, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(95, 3).Value = "Coquimbo"
$ws.Cells.Item(95, 4).Value = 45147
$ws.Cells.Item(95, 4).NumberFormat = $ws.Cells.Item(96, 4).NumberFormat
$ws.Cells.Item(95, 5).Value = 4
$ws.Cells.Item(95, 6).Value = 100112022
$ws.Cells.Item(95, 7).Value = "Arveja Verde"
$ws.Cells.Item(95, 8).Value = "Perfection"
$ws.Cells.Item(95, 9).Value = "Primera"
$ws.Cells.Item(95, 10).Value = 700
$ws.Cells.Item(95, 11).Value = 20000
$ws.Cells.Item(95, 12).Value = 22000
$ws.Cells.Item(95, 13).Value = 21000
$ws.Cells.Item(95, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(95, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(95, 16).Value = 840
$ws.Cells.Item(95, 17).Value = 25
$ws.Cells.Item(95, 18).Value = "Hortaliza"
